$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: Right count 3 -> 5
$ws.Range("B11").Value = 5

# Total row: Right count 69 -> 115
$ws.Range("B12").Value = 115

# Total row: Correct/Total marks string "67/84" -> "115/140"
$ws.Range("E12").Value = "115/140"
